# The DB test-case row (row 2) stores its AddifyVariables payload in column
# P. Update the demo quoteId used by the "InsuranceQuoteByDB" case from the
# old value (32633) to the new one (184).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSS-Accept-DB")

$ws.Range("P2").Value = "quoteId=184"

# Leave the sheet with that cell selected/active, matching where the edit
# was made.
$ws.Range("P2").Select()
